$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Status (B6) and Date (B8) values ---
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B6").Value = "draft"
$metaSheet.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Re-apply wrap-text alignment so the header/body styles record
#     an explicit "applied alignment" (wrapText + vertical top), matching
#     the already-visible vertical=top/wrapText formatting used throughout
#     both worksheets. ---
$metaHeader = $metaSheet.Range("A1:B1")
$metaHeader.WrapText = $true

$metaBody = $metaSheet.Range("A2:B21")
$metaBody.WrapText = $true

$conceptsSheet = $wb.Worksheets.Item("Concepts")
$conceptsHeader = $conceptsSheet.Range("A1:D1")
$conceptsHeader.WrapText = $true

$conceptsBody = $conceptsSheet.Range("A2:D4")
$conceptsBody.WrapText = $true
